$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.460.37'
$ws.Range("E2").Value = '  +0.83%  '
$ws.Range("D3").Value = '1.879.01'
$ws.Range("E3").Value = '  +1.23%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7197'
$ws.Range("E5").Value = '  +2.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '240.39'
$ws.Range("E6").Value = '  +0.93%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9996'
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3128'
$ws.Range("E8").Value = '  +3.41%  '
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07822'
$ws.Range("E9").Value = '  -2.23%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '25.12'
$ws.Range("E10").Value = '  +7.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08244'
$ws.Range("E11").Value = '  +0.73%  '
$ws.Range("D12").Value = '1.883.90'
$ws.Range("E12").Value = '  +4.11%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.7293'
$ws.Range("E13").Value = '  +4.01%  '
$ws.Range("E14").Value = '  +2.36%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.39'
$ws.Range("E15").Value = '  +2.08%  '
$ws.Range("D16").Value = '29.483.56'
$ws.Range("E16").Value = '  +1.30%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.939'
$ws.Range("E17").Value = '  +2.47%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '248.84'
$ws.Range("E18").Value = '  +4.90%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007885'
$ws.Range("E19").Value = '  +0.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.31'
$ws.Range("E20").Value = '  +0.81%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9989'
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.953'
$ws.Range("E22").Value = '  +6.82%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9996'
$ws.Range("E23").Value = '  -0.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.1568'
$ws.Range("E24").Value = '  +9.52%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '164.00'
$ws.Range("E25").Value = '  +0.80%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.059'
$ws.Range("E26").Value = '  +1.93%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.35'
$ws.Range("E27").Value = '  +1.45%  '
$ws.Range("E28").Value = '  -3.69%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.490'
$ws.Range("E29").Value = '  +0.81%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.390'
$ws.Range("E30").Value = '  +0.80%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.155'
$ws.Range("E31").Value = '  +3.40%  '
$ws.Range("E32").Value = '  +1.94%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.948'
$ws.Range("E33").Value = '  +1.87%  '
$ws.Range("E34").Value = '  +3.98%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7233'
$ws.Range("E35").Value = '  +1.54%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.676'
$ws.Range("E36").Value = '  +0.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01866'
$ws.Range("E37").Value = '  +0.82%  '
$ws.Range("D38").Value = '1.234.71'
$ws.Range("E38").Value = '  +8.51%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.722'
$ws.Range("E39").Value = '  +0.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9080'
$ws.Range("E40").Value = '  -2.65%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '73.85'
$ws.Range("E41").Value = '  +5.36%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.099'
$ws.Range("E42").Value = '  +3.14%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9995'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '103.77'
$ws.Range("E44").Value = '  +1.36%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5335'
$ws.Range("E45").Value = '  +0.10%  '
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '2.023.57'
$ws.Range("E46").Value = '  +5.28%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.766'
$ws.Range("E47").Value = '  +0.44%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.00000000121'
$ws.Range("E48").Value = '  +2.58%  '
$ws.Range("B49").Value = 'SynthetixNetwork'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.911'
$ws.Range("E49").Value = '  +12.57%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.299'
$ws.Range("E50").Value = '  +1.48%  '
$ws.Range("B51").Value = 'TheSandbox'
$ws.Range("C51").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4334'
$ws.Range("E51").Value = '  +2.08%  '
